$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 1.503062666666667
$ws.Range("H2").Value = 4.509188
$ws.Range("I2").Value = 0.1529782459008696
$ws.Range("J2").Value = 0.1529782459008696
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.855085666666667
$ws.Range("N2").Value = 17.565257
$ws.Range("O2").Value = 0.106330777065451
$ws.Range("P2").Value = 0.106330777065451
$ws.Range("Q2").Value = 8.800560675701778
$ws.Range("R2").Value = 79.205046081316
$ws.Range("S2").Value = 0.01626629576074912
$ws.Range("T2").Value = 0.01626629576074911
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 1.503062666666667
$ws.Range("H3").Value = 4.509188
$ws.Range("I3").Value = 0.1529782459008696
$ws.Range("J3").Value = 0.1529782459008696
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 33.421606
$ws.Range("N3").Value = 100.264818
$ws.Range("O3").Value = 0.6069501863972739
$ws.Range("P3").Value = 0.6069501863972738
$ws.Range("Q3").Value = 50.23476823864266
$ws.Range("R3").Value = 452.112914147784
$ws.Range("S3").Value = 0.09285017486426082
$ws.Range("T3").Value = 0.09285017486426081
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 1.503062666666667
$ws.Range("H4").Value = 4.509188
$ws.Range("I4").Value = 0.1529782459008696
$ws.Range("J4").Value = 0.1529782459008696
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.5207056666666666
$ws.Range("N4").Value = 1.562117
$ws.Range("O4").Value = 0.009456230243437438
$ws.Range("P4").Value = 0.009456230243437436
$ws.Range("Q4").Value = 0.7826532478884444
$ws.Range("R4").Value = 7.043879230995999
$ws.Range("S4").Value = 0.001446597515475813
$ws.Range("T4").Value = 0.001446597515475812
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 1.503062666666667
$ws.Range("H5").Value = 4.509188
$ws.Range("I5").Value = 0.1529782459008696
$ws.Range("J5").Value = 0.1529782459008696
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 15.267428
$ws.Range("N5").Value = 45.802284
$ws.Range("O5").Value = 0.2772628062938376
$ws.Range("P5").Value = 0.2772628062938375
$ws.Range("Q5").Value = 22.94790104282134
$ws.Range("R5").Value = 206.531109385392
$ws.Range("S5").Value = 0.04241517776038387
$ws.Range("T5").Value = 0.04241517776038386
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.420618666666667
$ws.Range("H6").Value = 4.261856
$ws.Range("I6").Value = 0.1445872860395478
$ws.Range("J6").Value = 0.1445872860395478
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.855085666666667
$ws.Range("N6").Value = 17.565257
$ws.Range("O6").Value = 0.106330777065451
$ws.Range("P6").Value = 0.106330777065451
$ws.Range("Q6").Value = 8.31784399299911
$ws.Range("R6").Value = 74.860595936992
$ws.Range("S6").Value = 0.01537407847836976
$ws.Range("T6").Value = 0.01537407847836976
# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.420618666666667
$ws.Range("H7").Value = 4.261856
$ws.Range("I7").Value = 0.1445872860395478
$ws.Range("J7").Value = 0.1445872860395478
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 33.421606
$ws.Range("N7").Value = 100.264818
$ws.Range("O7").Value = 0.6069501863972739
$ws.Range("P7").Value = 0.6069501863972738
$ws.Range("Q7").Value = 47.47935735357866
$ws.Range("R7").Value = 427.3142161822079
$ws.Range("S7").Value = 0.08775728021237951
$ws.Range("T7").Value = 0.0877572802123795
# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.420618666666667
$ws.Range("H8").Value = 4.261856
$ws.Range("I8").Value = 0.1445872860395478
$ws.Range("J8").Value = 0.1445872860395478
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.5207056666666666
$ws.Range("N8").Value = 1.562117
$ws.Range("O8").Value = 0.009456230243437438
$ws.Range("P8").Value = 0.009456230243437436
$ws.Range("Q8").Value = 0.7397241899057777
$ws.Range("R8").Value = 6.657517709151999
$ws.Range("S8").Value = 0.001367250667063712
$ws.Range("T8").Value = 0.001367250667063712
# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.420618666666667
$ws.Range("H9").Value = 4.261856
$ws.Range("I9").Value = 0.1445872860395478
$ws.Range("J9").Value = 0.1445872860395478
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 15.267428
$ws.Range("N9").Value = 45.802284
$ws.Range("O9").Value = 0.2772628062938376
$ws.Range("P9").Value = 0.2772628062938375
$ws.Range("Q9").Value = 21.68919320878933
$ws.Range("R9").Value = 195.202738879104
$ws.Range("S9").Value = 0.04008867668173483
$ws.Range("T9").Value = 0.04008867668173483
# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.029369333333333
$ws.Range("H10").Value = 12.088108
$ws.Range("I10").Value = 0.410099902266277
$ws.Range("J10").Value = 0.410099902266277
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.855085666666667
$ws.Range("N10").Value = 17.565257
$ws.Range("O10").Value = 0.106330777065451
$ws.Range("P10").Value = 0.106330777065451
$ws.Range("Q10").Value = 23.59230262930622
$ws.Range("R10").Value = 212.330723663756
$ws.Range("S10").Value = 0.04360624128243876
$ws.Range("T10").Value = 0.04360624128243875
# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 4.029369333333333
$ws.Range("H11").Value = 12.088108
$ws.Range("I11").Value = 0.410099902266277
$ws.Range("J11").Value = 0.410099902266277
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 33.421606
$ws.Range("N11").Value = 100.264818
$ws.Range("O11").Value = 0.6069501863972739
$ws.Range("P11").Value = 0.6069501863972738
$ws.Range("Q11").Value = 134.6679942871493
$ws.Range("R11").Value = 1212.011948584344
$ws.Range("S11").Value = 0.2489102121220207
$ws.Range("T11").Value = 0.2489102121220206
# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 4.029369333333333
$ws.Range("H12").Value = 12.088108
$ws.Range("I12").Value = 0.410099902266277
$ws.Range("J12").Value = 0.410099902266277
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.5207056666666666
$ws.Range("N12").Value = 1.562117
$ws.Range("O12").Value = 0.009456230243437438
$ws.Range("P12").Value = 0.009456230243437436
$ws.Range("Q12").Value = 2.098115444959555
$ws.Range("R12").Value = 18.883039004636
$ws.Range("S12").Value = 0.003877999098641106
$ws.Range("T12").Value = 0.003877999098641106
# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 4.029369333333333
$ws.Range("H13").Value = 12.088108
$ws.Range("I13").Value = 0.410099902266277
$ws.Range("J13").Value = 0.410099902266277
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 15.267428
$ws.Range("N13").Value = 45.802284
$ws.Range("O13").Value = 0.2772628062938376
$ws.Range("P13").Value = 0.2772628062938375
$ws.Range("Q13").Value = 61.51810618207467
$ws.Range("R13").Value = 553.662955638672
$ws.Range("S13").Value = 0.1137054497631765
$ws.Range("T13").Value = 0.1137054497631765
# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 2.872285333333334
$ws.Range("H14").Value = 8.616856
$ws.Range("I14").Value = 0.2923345657933055
$ws.Range("J14").Value = 0.2923345657933055
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 5.855085666666667
$ws.Range("N14").Value = 17.565257
$ws.Range("O14").Value = 0.106330777065451
$ws.Range("P14").Value = 0.106330777065451
$ws.Range("Q14").Value = 16.81747668577689
$ws.Range("R14").Value = 151.357290171992
$ws.Range("S14").Value = 0.0310841615438934
$ws.Range("T14").Value = 0.03108416154389339
# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 2.872285333333334
$ws.Range("H15").Value = 8.616856
$ws.Range("I15").Value = 0.2923345657933055
$ws.Range("J15").Value = 0.2923345657933055
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 33.421606
$ws.Range("N15").Value = 100.264818
$ws.Range("O15").Value = 0.6069501863972739
$ws.Range("P15").Value = 0.6069501863972738
$ws.Range("Q15").Value = 95.99638873024533
$ws.Range("R15").Value = 863.967498572208
$ws.Range("S15").Value = 0.1774325191986129
$ws.Range("T15").Value = 0.1774325191986129
# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 2.872285333333334
$ws.Range("H16").Value = 8.616856
$ws.Range("I16").Value = 0.2923345657933055
$ws.Range("J16").Value = 0.2923345657933055
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.5207056666666666
$ws.Range("N16").Value = 1.562117
$ws.Range("O16").Value = 0.009456230243437438
$ws.Range("P16").Value = 0.009456230243437436
$ws.Range("Q16").Value = 1.495615249350222
$ws.Range("R16").Value = 13.460537244152
$ws.Range("S16").Value = 0.002764382962256807
$ws.Range("T16").Value = 0.002764382962256807
# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 2.872285333333334
$ws.Range("H17").Value = 8.616856
$ws.Range("I17").Value = 0.2923345657933055
$ws.Range("J17").Value = 0.2923345657933055
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 15.267428
$ws.Range("N17").Value = 45.802284
$ws.Range("O17").Value = 0.2772628062938376
$ws.Range("P17").Value = 0.2772628062938375
$ws.Range("Q17").Value = 43.85240952212267
$ws.Range("R17").Value = 394.671685699104
$ws.Range("S17").Value = 0.0810535020885424
$ws.Range("T17").Value = 0.08105350208854238
